{"js": "const replacements = [\n  [\"171\u00d75=\", \"644\u00d75=\"],\n  [\"161\u00d73=\", \"234\u00d76=\"],\n  [\"935\u00d79=\", \"284\u00d72=\"],\n  [\"964\u00d78=\", \"766\u00d74=\"],\n  [\"616\u00d75=\", \"437\u00d78=\"],\n  [\"157\u00d75=\", \"404\u00d73=\"],\n  [\"431\u00d78=\", \"463\u00d73=\"],\n  [\"133\u00d74=\", \"396\u00d79=\"],\n  [\"479\u00d79=\", \"355\u00d76=\"],\n  [\"690\u00d79=\", \"625\u00d73=\"],\n  [\"253\u00d73=\", \"839\u00d73=\"],\n  [\"365\u00d72=\", \"194\u00d77=\"],\n  [\"668\u00d77=\", \"254\u00d78=\"],\n  [\"764\u00d76=\", \"458\u00d74=\"],\n  [\"684\u00d79=\", \"410\u00d74=\"],\n  [\"742\u00d72=\", \"576\u00d78=\"],\n  [\"156\u00d75=\", \"472\u00d78=\"],\n  [\"475\u00d77=\", \"192\u00d74=\"],\n  [\"966\u00d72=\", \"205\u00d79=\"],\n  [\"933\u00d78=\", \"995\u00d79=\"],\n  [\"904\u00d78=\", \"795\u00d74=\"],\n  [\"929\u00d77=\", \"259\u00d77=\"],\n  [\"634\u00d79=\", \"890\u00d78=\"],\n  [\"651\u00d79=\", \"507\u00d79=\"],\n  [\"249\u00d79=\", \"781\u00d72=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"171\u00d75=\", \"644\u00d75=\"),\n    @(\"161\u00d73=\", \"234\u00d76=\"),\n    @(\"935\u00d79=\", \"284\u00d72=\"),\n    @(\"964\u00d78=\", \"766\u00d74=\"),\n    @(\"616\u00d75=\", \"437\u00d78=\"),\n    @(\"157\u00d75=\", \"404\u00d73=\"),\n    @(\"431\u00d78=\", \"463\u00d73=\"),\n    @(\"133\u00d74=\", \"396\u00d79=\"),\n    @(\"479\u00d79=\", \"355\u00d76=\"),\n    @(\"690\u00d79=\", \"625\u00d73=\"),\n    @(\"253\u00d73=\", \"839\u00d73=\"),\n    @(\"365\u00d72=\", \"194\u00d77=\"),\n    @(\"668\u00d77=\", \"254\u00d78=\"),\n    @(\"764\u00d76=\", \"458\u00d74=\"),\n    @(\"684\u00d79=\", \"410\u00d74=\"),\n    @(\"742\u00d72=\", \"576\u00d78=\"),\n    @(\"156\u00d75=\", \"472\u00d78=\"),\n    @(\"475\u00d77=\", \"192\u00d74=\"),\n    @(\"966\u00d72=\", \"205\u00d79=\"),\n    @(\"933\u00d78=\", \"995\u00d79=\"),\n    @(\"904\u00d78=\", \"795\u00d74=\"),\n    @(\"929\u00d77=\", \"259\u00d77=\"),\n    @(\"634\u00d79=\", \"890\u00d78=\"),\n    @(\"651\u00d79=\", \"507\u00d79=\"),\n    @(\"249\u00d79=\", \"781\u00d72=\"),\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($pair[0], $false, $false, $false, $false, $false, $true, 1, $false, $pair[1], 2)\n}\n"}
